# TC10_Canine_StudyUBC-Breed_Diagnosis_PrimDiseaseSite.xlsx
# "Timing issue fix - keywords, updated tc1,2 in ubc01"
#
# Semantic edit: the CasesTab Neo4j query stored in cell B2 drops its
# trailing "Cohort" column (the OPTIONAL MATCH (co:cohort) clause stays,
# but the cohort_description RETURN expression and its label are removed).
# The sheet's viewport/selection also moved back up to B2, and the row
# height for row 2 shrinks to match the now-shorter text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the CasesTab query text in B2: drop the trailing Cohort column ---
$b2 = $ws.Range("B2")
$oldText = $b2.Value2
$cohortSuffix = ",`n        coalesce(co.cohort_description, '') AS ``Cohort``"
if ($oldText.Contains($cohortSuffix)) {
    $newText = $oldText.Replace($cohortSuffix, "")
    $b2.Value = $newText
}

# --- 2. Row 2 height shrinks (304.5 -> 290) now that the text is one line shorter ---
$ws.Rows.Item(2).RowHeight = 290

# --- 3. Viewport / selection moves from C8 back to B2 ---
$b2.Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
